$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
# B2 stays a text value ("2") rather than a number, matching the source data
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = -0.212
$ws.Range("E2").Value = 0.18
$ws.Range("G2").Value = 3.828823529411765
$ws.Range("H2").Value = 3.828823529411765
$ws.Range("I2").Value = 0.8162696868917794
$ws.Range("J2").Value = 0.7268745927835407
$ws.Range("K2").Value = 110.27
$ws.Range("L2").Value = 0.720718954248366
$ws.Range("M2").Value = 15.014
$ws.Range("N2").Value = 0.01671751475336822
$ws.Range("O2").Value = 0.1361567062664369
$ws.Range("P2").Value = 10.874
$ws.Range("Q2").Value = 0.0121077830976506
$ws.Range("R2").Value = 0.09861249659925637
$ws.Range("S2").Value = 4.140000000000001
$ws.Range("T2").Value = 0.2757426402024777
$ws.Range("U2").Value = 67.1
$ws.Range("V2").Value = 0.07471328359870838
$ws.Range("W2").Value = 0.9698853020433085
$ws.Range("X2").Value = 0.05095271576926303
$ws.Range("Y2").Value = 0.9189325862740454
$ws.Range("Z2").Value = 0.08591554036338774
$ws.Range("AA2").Value = -5.792663987835524
$ws.Range("AB2").Value = 0.04761096895172758
$ws.Range("AC2").Value = -5.840274956787251
$ws.Range("AD2").Value = 234.15
$ws.Range("AE2").Value = 1.673689527788807
$ws.Range("AF2").Value = 235.8236895277888
$ws.Range("AG2").Value = 168.7236895277888
$ws.Range("AH2").Value = 0.2079713932301699
$ws.Range("AI2").Value = 0.1271476128003893
$ws.Range("AJ2").Value = 0.1581551770775464
$ws.Range("AK2").Value = 0.09438434415263214
$ws.Range("AL2").Value = 15.8
$ws.Range("AM2").Value = 15.025
$ws.Range("AN2").Value = 1.86671875249135
$ws.Range("AO2").Value = 7.910126582278481
$ws.Range("AP2").Value = 1.345119262144146
$ws.Range("AQ2").Value = 8.318136439267887

# Row 3 updates
$ws.Range("D3").Value = -0.212
$ws.Range("E3").Value = 0.18
$ws.Range("G3").Value = 4.273062730627307
$ws.Range("H3").Value = 4.273062730627307
$ws.Range("I3").Value = 0.8797047970479706
$ws.Range("J3").Value = 0.8797047970479706
$ws.Range("K3").Value = 105.1
$ws.Range("L3").Value = 0.7756457564575645
$ws.Range("M3").Value = 14.14
$ws.Range("N3").Value = 0.02125356981812716
$ws.Range("O3").Value = 0.1345385347288297
$ws.Range("P3").Value = 10
$ws.Range("Q3").Value = 0.01503081316699234
$ws.Range("R3").Value = 0.09514747859181732
$ws.Range("S3").Value = 4.140000000000001
$ws.Range("T3").Value = 0.2927864214992928
$ws.Range("U3").Value = 47.2
$ws.Range("V3").Value = 0.07094543814820382
$ws.Range("W3").Value = 0.06658219828951537
$ws.Range("X3").Value = 0.05610132186345815
$ws.Range("Y3").Value = 0.01048087642605722
$ws.Range("Z3").Value = 0.07607231080170672
$ws.Range("AA3").Value = 0.06692117673478554
$ws.Range("AB3").Value = 0.04950008391969309
$ws.Range("AC3").Value = 0.01742109281509245
$ws.Range("AD3").Value = 232.6
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 232.6
$ws.Range("AG3").Value = 185.4
$ws.Range("AH3").Value = 0.2590488918587816
$ws.Range("AI3").Value = 0.1270413457862254
$ws.Range("AJ3").Value = 0.2179381685670624
$ws.Range("AK3").Value = 0.103941245725178
$ws.Range("AL3").Value = 15.8
$ws.Range("AM3").Value = 15.8
$ws.Range("AN3").Value = 1.949706621961442
$ws.Range("AO3").Value = 7.544303797468355
$ws.Range("AP3").Value = 1.55406538139145
$ws.Range("AQ3").Value = 7.544303797468355

# Row 4 new row
$ws.Range("A4").Value = "Indonesia"
$ws.Range("B4").Value = "PT Ashmore Asset Management Indonesia Tbk (IDX:AMOR)"
$ws.Range("C4").Value = "Investments & Asset Management"
$ws.Range("G4").Value = 0.3891428571428571
$ws.Range("H4").Value = 0.3891428571428571
$ws.Range("I4").Value = 0.3251006911109851
$ws.Range("J4").Value = 0.2538928358072194
$ws.Range("K4").Value = 5.17
$ws.Range("L4").Value = 0.2954285714285714
$ws.Range("M4").Value = 0.874
$ws.Range("N4").Value = 0.003754295532646048
$ws.Range("O4").Value = 0.1690522243713733
$ws.Range("P4").Value = 0.874
$ws.Range("Q4").Value = 0.003754295532646048
$ws.Range("R4").Value = 0.1690522243713733
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 19.9
$ws.Range("V4").Value = 0.08548109965635738
$ws.Range("W4").Value = 1.873188405797102
$ws.Range("X4").Value = 0.04580410967506791
$ws.Range("Y4").Value = 1.827384296122034
$ws.Range("Z4").Value = -45.89435978120066
$ws.Range("AA4").Value = -11.65224915240583
$ws.Range("AB4").Value = 0.04572185398376206
$ws.Range("AC4").Value = -11.6979710063896
$ws.Range("AD4").Value = 1.55
$ws.Range("AE4").Value = 1.673689527788807
$ws.Range("AF4").Value = 3.223689527788808
$ws.Range("AG4").Value = -16.67631047221119
$ws.Range("AH4").Value = 0.01365833037454174
$ws.Range("AI4").Value = 0.1353144534572859
$ws.Range("AJ4").Value = -0.07716095587969768
$ws.Range("AK4").Value = -4.250160557838302
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = -0.775
$ws.Range("AN4").Value = 0.2526899250081513
$ws.Range("AP4").Value = -2.718668156539158
$ws.Range("AQ4").Value = -7.458064516129032
